$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'schubert-winterreise_0'
$ws.Range('B2').Value = 'jaah_87'
$ws.Range("C2").Value = 0.07051282051282051
$ws.Range('D2').Value = '[[''B:min'', ''F#:maj'', ''B:min'', ''B:7'']]'
$ws.Range('E2').Value = '[[''F:min'', ''C'', ''F:min'', ''F:7'']]'
$ws.Range('F2').Value = '[(12.14, 21.58)]'
$ws.Range('G2').Value = '[(11.89, 19.531)]'
$ws.Range('H2').Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Range("I2").Value = ""

# Row 3
$ws.Range('A3').Value = 'schubert-winterreise_187'
$ws.Range('B3').Value = 'schubert-winterreise_10'
$ws.Range("C3").Value = 0.08749999999999999
$ws.Range('D3').Value = '[[''C:min'', ''F:7/C'', ''A#''], [''C/G'', ''G:7'', ''C'']]'
$ws.Range('E3').Value = '[[''C:min'', ''F:7/A'', ''A#:maj''], [''C:maj'', ''G:7'', ''C:maj'']]'
$ws.Range('F3').Value = '[(37.2, 39.9), (236.4, 239.98)]'
$ws.Range('G3').Value = '[(37.2, 43.88), (0.8, 9.28)]'
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

# Row 4
$ws.Range('A4').Value = 'schubert-winterreise_205'
$ws.Range('B4').Value = 'schubert-winterreise_17'
$ws.Range("C4").Value = 0.323076923076923
$ws.Range('D4').Value = '[[''G:maj'', ''C:maj/G'', ''G:maj'']]'
$ws.Range('E4').Value = '[[''D:maj/A'', ''G:maj'', ''D:maj/A'']]'
$ws.Range('F4').Value = '[(16.92, 23.9)]'
$ws.Range('G4').Value = '[(138.02, 142.34)]'
$ws.Range('H4').Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I4").Value = ""

# Row 5
$ws.Range('A5').Value = 'schubert-winterreise_182'
$ws.Range('B5').Value = 'schubert-winterreise_145'
$ws.Range("C5").Value = 1
$ws.Range('D5').Value = '[[''G:maj'', ''E:7/G#'', ''A:min'', ''D:7/F#'', ''G:maj'']]'
$ws.Range('E5').Value = '[[''G:maj'', ''E:7/G#'', ''A:min'', ''D:7/F#'', ''G:maj'']]'
$ws.Range('F5').Value = '[(82.24, 89.66)]'
$ws.Range('G5').Value = '[(76.74, 83.98)]'
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

# Row 6
$ws.Range('A6').Value = 'schubert-winterreise_11'
$ws.Range('B6').Value = 'isophonics_274'
$ws.Range("C6").Value = 0.1083333333333333
$ws.Range('D6').Value = '[[''C:maj'', ''F:maj'', ''C:maj/E'']]'
$ws.Range('E6').Value = '[[''Ab'', ''Db'', ''Ab'']]'
$ws.Range('F6').Value = '[(69.86, 73.74)]'
$ws.Range('G6').Value = '[(3.129454, 5.172811)]'
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

# Row 7
$ws.Range('A7').Value = 'isophonics_6'
$ws.Range('B7').Value = 'isophonics_277'
$ws.Range("C7").Value = 0.2095238095238095
$ws.Range('D7').Value = '[[''G/5'', ''D'', ''G/5'', ''D'']]'
$ws.Range('E7').Value = '[[''A'', ''E'', ''A'', ''E'']]'
$ws.Range('F7').Value = '[(21.391, 26.25)]'
$ws.Range('G7').Value = '[(31.840929, 39.410634)]'
$ws.Range("H7").Value = ""
$ws.Range('I7').Value = 'spotify:track:2RnPATK99oGOZygnD2GTO6'

# Row 8
$ws.Range('A8').Value = 'schubert-winterreise_116'
$ws.Range('B8').Value = 'schubert-winterreise_179'
$ws.Range("C8").Value = 0.3342175066312997
$ws.Range('D8').Value = '[[''D:maj/G'', ''G:min'', ''D:maj/G'', ''G:min'', ''D:maj/G'', ''G:min'']]'
$ws.Range('E8').Value = '[[''A:maj'', ''D:min'', ''A:maj'', ''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range('F8').Value = '[(44.3, 67.1)]'
$ws.Range('G8').Value = '[(0.22, 8.88)]'
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""

# Row 9
$ws.Range('A9').Value = 'schubert-winterreise_65'
$ws.Range('B9').Value = 'schubert-winterreise_105'
$ws.Range("C9").Value = 0.2528735632183908
$ws.Range('D9').Value = '[[''D#:min'', ''A#:7'', ''D#:min'', ''A#:maj'', ''D#:min'']]'
$ws.Range('E9').Value = '[[''G:min'', ''D:7/G'', ''G:min'', ''D:maj/G'', ''G:min'']]'
$ws.Range('F9').Value = '[(12.32, 19.9)]'
$ws.Range('G9').Value = '[(23.74, 51.72)]'
$ws.Range('H9').Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I9").Value = ""

# Row 10
$ws.Range('A10').Value = 'jaah_21'
$ws.Range('B10').Value = 'jaah_28'
$ws.Range("C10").Value = 0.111038961038961
$ws.Range('D10').Value = '[[''Eb:7'', ''F:min7'', ''F:min7'']]'
$ws.Range('E10').Value = '[[''Ab:7'', ''Db:maj6'', ''Db:maj6'']]'
$ws.Range('F10').Value = '[(103.22, 107.52)]'
$ws.Range('G10').Value = '[(16.96, 19.46)]'
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""

# Row 11
$ws.Range('A11').Value = 'schubert-winterreise_184'
$ws.Range('B11').Value = 'schubert-winterreise_6'
$ws.Range("C11").Value = 0.07334525939177103
$ws.Range('D11').Value = '[[''A#:maj/F'', ''F:7'', ''A#:maj'']]'
$ws.Range('E11').Value = '[[''D:maj/F#'', ''A:7'', ''D:maj'']]'
$ws.Range('F11').Value = '[(28.02, 30.1)]'
$ws.Range('G11').Value = '[(37.54, 42.72)]'
$ws.Range('H11').Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'
$ws.Range('I11').Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 12
$ws.Range('A12').Value = 'jaah_21'
$ws.Range('B12').Value = 'isophonics_204'
$ws.Range("C12").Value = 0.1675824175824176
$ws.Range('D12').Value = '[[''Ab/b5'', ''Ab:7'', ''Db/3'', ''Db:min/b3'', ''Ab'']]'
$ws.Range('E12').Value = '[[''A'', ''A:7'', ''D'', ''D:min'', ''A'']]'
$ws.Range('F12').Value = '[(5.32, 11.34)]'
$ws.Range('G12').Value = '[(18.651995, 28.810725)]'
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

# Row 13
$ws.Range('A13').Value = 'isophonics_224'
$ws.Range('B13').Value = 'isophonics_201'
$ws.Range("C13").Value = 0.1302083333333333
$ws.Range('D13').Value = '[[''E:7'', ''A'', ''D'']]'
$ws.Range('E13').Value = '[[''C:7'', ''F:maj'', ''Bb:maj'']]'
$ws.Range('F13').Value = '[(26.987959, 33.059977)]'
$ws.Range('G13').Value = '[(63.251, 70.247)]'
$ws.Range('H13').Value = 'spotify:track:6dGnYIeXmHdcikdzNNDMm2'
$ws.Range("I13").Value = ""

# Row 14
$ws.Range('A14').Value = 'schubert-winterreise_211'
$ws.Range('B14').Value = 'schubert-winterreise_63'
$ws.Range("C14").Value = 0.233974358974359
$ws.Range('D14').Value = '[[''F:maj'', ''C:7/E'', ''F:maj'', ''C:7/E'', ''F:maj'']]'
$ws.Range('E14').Value = '[[''D#/G'', ''A#:7/F'', ''D#/G'', ''A#:7'', ''D#'']]'
$ws.Range('F14').Value = '[(34.64, 41.66)]'
$ws.Range('G14').Value = '[(43.32, 47.42)]'
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = ""

# Row 15
$ws.Range('A15').Value = 'isophonics_288'
$ws.Range('B15').Value = 'isophonics_294'
$ws.Range("C15").Value = 0.08452380952380953
$ws.Range('D15').Value = '[[''E'', ''D'', ''A'']]'
$ws.Range('E15').Value = '[[''G'', ''F'', ''C/3'']]'
$ws.Range('F15').Value = '[(22.169818, 24.921383)]'
$ws.Range('G15').Value = '[(45.670113, 48.251979)]'
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# Row 16
$ws.Range('A16').Value = 'schubert-winterreise_84'
$ws.Range('B16').Value = 'jaah_20'
$ws.Range("C16").Value = 0.1619047619047619
$ws.Range('D16').Value = '[[''F:min/C'', ''C:7'', ''F:min'', ''G#:maj'']]'
$ws.Range('E16').Value = '[[''F:min'', ''C:7'', ''F:min'', ''Ab'']]'
$ws.Range('F16').Value = '[(47.4, 56.94)]'
$ws.Range('G16').Value = '[(45.67, 51.02)]'
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""

# Row 17
$ws.Range('A17').Value = 'isophonics_194'
$ws.Range('B17').Value = 'isophonics_93'
$ws.Range("C17").Value = 0.06227967097532315
$ws.Range('D17').Value = '[[''G'', ''D'', ''D''], [''D'', ''A'', ''B:min''], [''G:maj6/5'', ''G/5'', ''D'']]'
$ws.Range('E17').Value = '[[''F:maj/9'', ''C'', ''C/7''], [''C'', ''G'', ''A:min''], [''D:min7/4'', ''F:maj/9'', ''C'']]'
$ws.Range('F17').Value = '[(10.919344, 16.143834), (60.888687, 69.665829), (4.893766, 6.623653)]'
$ws.Range('G17').Value = '[(20.323832, 24.050634), (3.059795, 10.455351), (19.081564, 22.808367)]'
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
